$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6956945
$ws.Range("I18").Value = 27777776
$ws.Range("K18").Value = 27777776
$ws.Range("M18").Value = -27777492

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3210
$ws.Range("I102").Value = 2700
$ws.Range("J102").Value = 5760
$ws.Range("K102").Value = 2700
$ws.Range("L102").Value = 5760
$ws.Range("M102").Value = -1078
$ws.Range("N102").Value = -9004

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1604.3636
$ws.Range("I20").Value = 1092.5238
$ws.Range("J20").Value = 2500.0833
$ws.Range("K20").Value = 1092.5238
$ws.Range("L20").Value = 2500.0833
$ws.Range("M20").Value = -845.5237999999999
$ws.Range("N20").Value = -2994.0833

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2234.5454
$ws.Range("I86").Value = 1506.1538
$ws.Range("J86").Value = 3286.6667
$ws.Range("K86").Value = 1506.1538
$ws.Range("L86").Value = 3286.6667
$ws.Range("M86").Value = -383.1538
$ws.Range("N86").Value = -5532.6667

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2234.5454
$ws.Range("I89").Value = 1506.1538
$ws.Range("J89").Value = 3286.6667
$ws.Range("K89").Value = 7530.769
$ws.Range("L89").Value = 16433.3335
$ws.Range("M89").Value = -1914.769
$ws.Range("N89").Value = -27665.3335

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4840
$ws.Range("I99").Value = 1100
$ws.Range("J99").Value = 7333.3335
$ws.Range("K99").Value = 1100
$ws.Range("L99").Value = 7333.3335
$ws.Range("M99").Value = 398
$ws.Range("N99").Value = -10329.3335

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4840
$ws.Range("I126").Value = 1100
$ws.Range("J126").Value = 7333.3335
$ws.Range("K126").Value = 3300
$ws.Range("L126").Value = 22000.0005
$ws.Range("M126").Value = -830
$ws.Range("N126").Value = -26940.0005

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4689
$ws.Range("I132").Value = 2870.6667
$ws.Range("J132").Value = 5598.1665
$ws.Range("K132").Value = 8612.000100000001
$ws.Range("L132").Value = 16794.4995
$ws.Range("M132").Value = -6082.000100000001
$ws.Range("N132").Value = -21854.4995

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 12219.637
$ws.Range("I64").Value = 1680.4
$ws.Range("J64").Value = 21002.334
$ws.Range("K64").Value = 5041.200000000001
$ws.Range("L64").Value = 63007.00199999999
$ws.Range("M64").Value = -4771.200000000001
$ws.Range("N64").Value = -63547.00199999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 12219.637
$ws.Range("I67").Value = 1680.4
$ws.Range("J67").Value = 21002.334
$ws.Range("K67").Value = 5041.200000000001
$ws.Range("L67").Value = 63007.00199999999
$ws.Range("M67").Value = -4105.200000000001
$ws.Range("N67").Value = -64879.00199999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 17497.857
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 19914.166
$ws.Range("K74").Value = 9000
$ws.Range("L74").Value = 59742.49800000001
$ws.Range("M74").Value = -7939
$ws.Range("N74").Value = -61864.49800000001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2254.5334
$ws.Range("I75").Value = 1141.6666
$ws.Range("J75").Value = 2532.75
$ws.Range("K75").Value = 3424.9998
$ws.Range("L75").Value = 7598.25
$ws.Range("M75").Value = -2426.9998
$ws.Range("N75").Value = -9594.25

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3647.5
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 3777
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 11331
$ws.Range("M76").Value = -8617
$ws.Range("N76").Value = -12097

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 17497.857
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 19914.166
$ws.Range("K77").Value = 27000
$ws.Range("L77").Value = 179227.494
$ws.Range("M77").Value = -21696
$ws.Range("N77").Value = -189835.494

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2254.5334
$ws.Range("I78").Value = 1141.6666
$ws.Range("J78").Value = 2532.75
$ws.Range("K78").Value = 10274.9994
$ws.Range("L78").Value = 22794.75
$ws.Range("M78").Value = -5282.999400000001
$ws.Range("N78").Value = -32778.75

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 3647.5
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 3777
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 11331
$ws.Range("M79").Value = -7674
$ws.Range("N79").Value = -13983

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 14375
$ws.Range("I81").Value = 750
$ws.Range("J81").Value = 28000
$ws.Range("K81").Value = 2250
$ws.Range("L81").Value = 84000
$ws.Range("M81").Value = -1127
$ws.Range("N81").Value = -86246

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 14375
$ws.Range("I84").Value = 750
$ws.Range("J84").Value = 28000
$ws.Range("K84").Value = 6750
$ws.Range("L84").Value = 252000
$ws.Range("M84").Value = -1134
$ws.Range("N84").Value = -263232

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10583.333
$ws.Range("I87").Value = 6920
$ws.Range("K87").Value = 20760
$ws.Range("M87").Value = -19512

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 10583.333
$ws.Range("I90").Value = 6920
$ws.Range("K90").Value = 62280
$ws.Range("M90").Value = -56040

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 5541
$ws.Range("I46").Value = 5541
$ws.Range("K46").Value = 5541
$ws.Range("M46").Value = -5385

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4640.2
$ws.Range("I113").Value = 4000.3333
$ws.Range("J113").Value = 5600
$ws.Range("K113").Value = 4000.3333
$ws.Range("L113").Value = 5600
$ws.Range("M113").Value = -1830.3333
$ws.Range("N113").Value = -9940

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2783.111
$ws.Range("I132").Value = 1833
$ws.Range("J132").Value = 3258.1667
$ws.Range("K132").Value = 5499
$ws.Range("L132").Value = 9774.500100000001
$ws.Range("M132").Value = -2969
$ws.Range("N132").Value = -14834.5001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 407131.56
$ws.Range("I132").Value = 627718.2
$ws.Range("J132").Value = 14977.556
$ws.Range("K132").Value = 1883154.6
$ws.Range("L132").Value = 44932.66800000001
$ws.Range("M132").Value = -1880624.6
$ws.Range("N132").Value = -49992.66800000001

Write-Host "Applied all cell updates"
